{"js": "// The document contains four `<id>p047v_aN</id>` markers, each split\n// across three separate runs: \"<id>\", \"p047v_aN\", \"</id>\".\n// The edit collapses each triple into a single run whose text is\n// \"<id>p047v_N</id>\" (the \"a\" before the digit is dropped), keeping\n// the formatting of the opening \"<id>\" run.\nconst body = context.document.body;\n\nfor (let i = 1; i <= 4; i++) {\n  const oldTag = \"<id>p047v_a\" + i + \"</id>\";\n  const newTag = \"<id>p047v_\" + i + \"</id>\";\n\n  const results = body.search(oldTag, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find \" + oldTag);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newTag, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains four `<id>p047v_aN</id>` markers, each split\n# across three separate runs: \"<id>\", \"p047v_aN\", \"</id>\".\n# The edit collapses each triple into a single run whose text is\n# \"<id>p047v_N</id>\" (the \"a\" before the digit is dropped), keeping\n# the formatting of the opening \"<id>\" run.\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le 4; $i++) {\n    $oldTag = \"<id>p047v_a$i</id>\"\n    $newTag = \"<id>p047v_$i</id>\"\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $oldTag\n    $find.MatchCase = $true\n    $found = $find.Execute()\n\n    if (-not $found) {\n        throw \"Could not find $oldTag\"\n    }\n\n    $range.Text = $newTag\n}\n"}
